# Refresh the crypto tracker sheet with the latest scraped price/volume(1h)
# figures, as produced by the 'Updated cryptos list' GitHub Actions run on
# Sat May 18 10:13:00 UTC 2024.
#
# Two coins also swapped rank position on the source site between scrapes:
#   rows 32-33: EthereumClassic <-> PEPE
#   rows 46-47: VeChain <-> Bittensor
# so those rows' Coin/Link/Price/Volume cells are rewritten wholesale rather
# than row-moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'67.170.18"
$ws.Range('E2').Value = '  +1.38%  '

$ws.Range('D3').Value = "'3.125.06"
$ws.Range('E3').Value = '  +3.36%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = "'580.15"
$ws.Range('E5').Value = '  +0.64%  '

$ws.Range('D6').Value = "'175.22"
$ws.Range('E6').Value = '  +4.40%  '

$ws.Range('D7').Value = "'1.00"

$ws.Range('D8').Value = "'3.120.40"
$ws.Range('E8').Value = '  +3.31%  '

$ws.Range('D9').Value = "'0.524"
$ws.Range('E9').Value = '  +1.12%  '

$ws.Range('D10').Value = "'6.48"
$ws.Range('E10').Value = '  -2.67%  '

$ws.Range('D11').Value = "'0.155"
$ws.Range('E11').Value = '  +2.01%  '

$ws.Range('D12').Value = "'0.484"
$ws.Range('E12').Value = '  +0.28%  '

$ws.Range('D13').Value = "'0.0000250"
$ws.Range('E13').Value = '  +1.22%  '

$ws.Range('D14').Value = "'37.31"
$ws.Range('E14').Value = '  +2.22%  '

$ws.Range('D15').Value = "'0.125"
$ws.Range('E15').Value = '  +0.29%  '

$ws.Range('D16').Value = "'3.640.89"
$ws.Range('E16').Value = '  +3.21%  '

$ws.Range('D17').Value = "'67.102.14"
$ws.Range('E17').Value = '  +1.25%  '

$ws.Range('D18').Value = "'7.19"
$ws.Range('E18').Value = '  -0.57%  '

$ws.Range('D19').Value = "'3.122.58"
$ws.Range('E19').Value = '  +3.17%  '

$ws.Range('D20').Value = "'16.33"
$ws.Range('E20').Value = '  -0.90%  '

$ws.Range('D21').Value = "'487.12"
$ws.Range('E21').Value = '  +4.26%  '

$ws.Range('E22').Value = '  +1.74%  '

$ws.Range('D23').Value = "'7.65"
$ws.Range('E23').Value = '  +3.76%  '

$ws.Range('D24').Value = "'84.43"
$ws.Range('E24').Value = '  +1.73%  '

$ws.Range('D25').Value = "'13.34"
$ws.Range('E25').Value = '  +4.22%  '

$ws.Range('D26').Value = "'2.34"
$ws.Range('E26').Value = '  +3.78%  '

$ws.Range('D27').Value = "'10.07"
$ws.Range('E27').Value = '  +0.22%  '

$ws.Range('E28').Value = '  -0.05%  '

$ws.Range('D29').Value = "'8.02"
$ws.Range('E29').Value = '  -2.36%  '

$ws.Range('E30').Value = '  -1.29%  '

$ws.Range('E31').Value = '  +2.20%  '

$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = "'0.0000101"
$ws.Range('E32').Value = '  +1.54%  '

$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = "'28.88"
$ws.Range('E33').Value = '  +2.29%  '

$ws.Range('E34').Value = '  -2.45%  '

$ws.Range('D36').Value = "'5.94"
$ws.Range('E36').Value = '  +1.30%  '

$ws.Range('D37').Value = "'0.992"
$ws.Range('E37').Value = '  +0.25%  '

$ws.Range('D38').Value = "'47.45"
$ws.Range('E38').Value = '  -1.53%  '

$ws.Range('D39').Value = "'2.12"
$ws.Range('E39').Value = '  +3.16%  '

$ws.Range('E40').Value = '  +1.55%  '

$ws.Range('D41').Value = "'0.315"
$ws.Range('E41').Value = '  +1.19%  '

$ws.Range('E42').Value = '  +2.16%  '

$ws.Range('D43').Value = "'8.68"
$ws.Range('E43').Value = '  +0.86%  '

$ws.Range('E44').Value = '  -0.72%  '

$ws.Range('D45').Value = "'2.853.12"
$ws.Range('E45').Value = '  +5.30%  '

$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = "'386.23"
$ws.Range('E46').Value = '  +1.82%  '

$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = "'0.0359"
$ws.Range('E47').Value = '  +0.07%  '

$ws.Range('D48').Value = "'136.57"
$ws.Range('E48').Value = '  +1.19%  '

$ws.Range('D50').Value = "'25.10"
$ws.Range('E50').Value = '  +2.66%  '

$ws.Range('D51').Value = "'2.23"
$ws.Range('E51').Value = '  +0.13%  '
